$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45177 -> 45178, i.e. 2023-09-08 -> 2023-09-09) for every data row (2-339).
$ws.Range("C2:C339").Value = 45178
